# Fruta / hortaliza, semanal
# Insert two new weekly price rows for "Mandarina" (Murcott) at the top of the
# data block (row 106), pushing the existing rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 106 (shifts rows 106:131 down to 108:133).
# Inserting this way copies formatting (e.g. the date number format on column D)
# from the row immediately above, just like Excel's normal UI behavior.
$ws.Rows.Item(106).Insert()
$ws.Rows.Item(107).Insert()

$fecha = Get-Date -Year 2021 -Month 9 -Day 24 -Hour 0 -Minute 0 -Second 0

# --- New row 106: Murcott / Primera ---
$ws.Range("A106").Value = 7
$ws.Range("B106").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C106").Value = "Ñuble"
$ws.Range("D106").Value = $fecha
$ws.Range("E106").Value = 16
$ws.Range("F106").Value = "Fruta"
$ws.Range("G106").Value = 100102
$ws.Range("H106").Value = "Cítricos"
$ws.Range("I106").Value = 100102004
$ws.Range("J106").Value = "Mandarina"
$ws.Range("K106").Value = "Murcott"
$ws.Range("L106").Value = "Primera"
$ws.Range("M106").Value = 240
$ws.Range("N106").Value = 5500
$ws.Range("O106").Value = 6000
$ws.Range("P106").Value = 5750
$ws.Range("Q106").Value = "$/bandeja 10 kilos"
$ws.Range("R106").Value = "Provincia de Limarí"
$ws.Range("S106").Value = 575
$ws.Range("T106").Value = 10

# --- New row 107: Murcott / Segunda ---
$ws.Range("A107").Value = 7
$ws.Range("B107").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C107").Value = "Ñuble"
$ws.Range("D107").Value = $fecha
$ws.Range("E107").Value = 16
$ws.Range("F107").Value = "Fruta"
$ws.Range("G107").Value = 100102
$ws.Range("H107").Value = "Cítricos"
$ws.Range("I107").Value = 100102004
$ws.Range("J107").Value = "Mandarina"
$ws.Range("K107").Value = "Murcott"
$ws.Range("L107").Value = "Segunda"
$ws.Range("M107").Value = 90
$ws.Range("N107").Value = 5000
$ws.Range("O107").Value = 5000
$ws.Range("P107").Value = 5000
$ws.Range("Q107").Value = "$/bandeja 10 kilos"
$ws.Range("R107").Value = "Provincia de Limarí"
$ws.Range("S107").Value = 500
$ws.Range("T107").Value = 10
